$d = $word.ActiveDocument

$replacements = @(
    @{old="589÷2=294, 1"; new="387÷2=193, 1"},
    @{old="568÷8=71, 0"; new="648÷2=324, 0"},
    @{old="412÷7=58, 6"; new="456÷5=91, 1"},
    @{old="792÷4=198, 0"; new="884÷3=294, 2"},
    @{old="711÷2=355, 1"; new="201÷2=100, 1"},
    @{old="281÷3=93, 2"; new="985÷5=197, 0"},
    @{old="736÷9=81, 7"; new="440÷3=146, 2"},
    @{old="105÷4=26, 1"; new="413÷7=59, 0"},
    @{old="118÷3=39, 1"; new="136÷5=27, 1"},
    @{old="507÷9=56, 3"; new="300÷7=42, 6"},
    @{old="698÷8=87, 2"; new="425÷2=212, 1"},
    @{old="827÷6=137, 5"; new="273÷6=45, 3"},
    @{old="444÷4=111, 0"; new="611÷8=76, 3"},
    @{old="238÷7=34, 0"; new="316÷7=45, 1"},
    @{old="221÷6=36, 5"; new="532÷3=177, 1"},
    @{old="904÷6=150, 4"; new="719÷4=179, 3"},
    @{old="908÷5=181, 3"; new="715÷7=102, 1"},
    @{old="467÷5=93, 2"; new="954÷5=190, 4"},
    @{old="277÷9=30, 7"; new="758÷5=151, 3"},
    @{old="774÷5=154, 4"; new="365÷8=45, 5"},
    @{old="645÷7=92, 1"; new="720÷2=360, 0"},
    @{old="842÷8=105, 2"; new="923÷4=230, 3"},
    @{old="526÷8=65, 6"; new="350÷5=70, 0"},
    @{old="580÷4=145, 0"; new="306÷8=38, 2"},
    @{old="440÷4=110, 0"; new="415÷5=83, 0"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
